$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 770 ("「エジプト特許庁」" entry), shifting subsequent rows up.
$ws.Rows.Item(770).Delete()
